$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 198, pushing existing rows 198-299 down to 199-300
$ws.Rows.Item(198).Insert()

# Populate the newly inserted row 198 with the new data record
$ws.Cells.Item(198, 1).Value = 6
$ws.Cells.Item(198, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(198, 3).Value = "Metropolitana"
$ws.Cells.Item(198, 4).Value = 44806
$ws.Cells.Item(198, 5).Value = 13
$ws.Cells.Item(198, 6).Value = 100112026
$ws.Cells.Item(198, 7).Value = "Haba"
$ws.Cells.Item(198, 8).Value = "Sin especificar"
$ws.Cells.Item(198, 9).Value = "Primera"
$ws.Cells.Item(198, 10).Value = 800
$ws.Cells.Item(198, 11).Value = 8000
$ws.Cells.Item(198, 12).Value = 10000
$ws.Cells.Item(198, 13).Value = 9125
$ws.Cells.Item(198, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(198, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(198, 16).Value = 365
$ws.Cells.Item(198, 17).Value = 25
$ws.Cells.Item(198, 18).Value = "Hortaliza"
